$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.602.57"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "'3.444.30"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'580.59"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "'175.16"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").Value = "'3.441.75"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").Value = "'6.86"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").Value = "'0.419"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").Value = "'4.038.62"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'30.96"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").Value = "'0.131"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "'66.555.32"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'0.0000172"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "'3.441.54"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -3.80%  "
$ws.Range("D20").Value = "'13.78"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").Value = "'375.20"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").Value = "'7.68"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'70.70"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "'9.82"
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'5.85"
$ws.Range("E31").Value = "  -5.09%  "
$ws.Range("D32").Value = "'23.87"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").Value = "'1.34"
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -4.71%  "
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("D38").Value = "'159.02"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'26.97"
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").Value = "'6.52"
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").Value = "'2.690.40"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("E46").Value = "  -4.26%  "
$ws.Range("D47").Value = "'25.11"
$ws.Range("E47").Value = "  -5.33%  "
$ws.Range("D48").Value = "'40.33"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").Value = "'320.00"
$ws.Range("E50").Value = "  -4.94%  "
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  -3.78%  "
